# Atualizacao dos dados BIBI (previsao_retorno) - atualizei dados bibi e add
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Simple "situacao" (column J) text refreshes: days-since-last-purchase
# figures nudged forward, no other column touched on these rows.
$ws.Range("J5").Value   = "INATIVO - 15.7 meses sem comprar"
$ws.Range("J6").Value   = "INATIVO - 17.2 meses sem comprar"
$ws.Range("J16").Value  = "INATIVO - 40.7 meses sem comprar"
$ws.Range("J24").Value  = "INATIVO - 38.7 meses sem comprar"
$ws.Range("J25").Value  = "INATIVO - 0.7 meses sem comprar"
$ws.Range("J31").Value  = "INATIVO - 7.7 meses sem comprar"
$ws.Range("J39").Value  = "INATIVO - 33.1 meses sem comprar"
$ws.Range("J46").Value  = "INATIVO - 7.0 meses sem comprar"
$ws.Range("J47").Value  = "INATIVO - 16.8 meses sem comprar"
$ws.Range("J51").Value  = "INATIVO - 8.4 meses sem comprar"
$ws.Range("J67").Value  = "INATIVO - 28.7 meses sem comprar"
$ws.Range("J71").Value  = "INATIVO - 12.2 meses sem comprar"
$ws.Range("J81").Value  = "INATIVO - 7.0 meses sem comprar"
$ws.Range("J90").Value  = "INATIVO - 15.9 meses sem comprar"
$ws.Range("J92").Value  = "INATIVO - 12.6 meses sem comprar"
$ws.Range("J93").Value  = "INATIVO - 12.0 meses sem comprar"
$ws.Range("J100").Value = "INATIVO - 33.8 meses sem comprar"
$ws.Range("J104").Value = "INATIVO - 38.1 meses sem comprar"
$ws.Range("J106").Value = "INATIVO - 15.4 meses sem comprar"

# --- Row 63 (id_cliente 6486 - CARLOS ALBERTO): new purchase recorded,
# probabilities / regularity / pattern / dates recalculated.
$ws.Range("B63").Value = 0.25
$ws.Range("C63").Value = 0.17
$ws.Range("E63").Value = 7
$ws.Range("G63").Value = "1x a cada 2 meses - irregular"
$ws.Range("H63").Value = (Get-Date -Year 2025 -Month 7 -Day 23 -Hour 10 -Minute 42 -Second 53)
$ws.Range("I63").Value = (Get-Date -Year 2025 -Month 9 -Day 23 -Hour 10 -Minute 42 -Second 53)

# --- Row 78 (id_cliente 20764 - EDILSON SOARES): client returned to
# "ATIVO" - total de compras, padrao, datas e situacao todos mudaram.
$ws.Range("E78").Value = 7
$ws.Range("G78").Value = "1x a cada 6 meses - irregular"
$ws.Range("H78").Value = (Get-Date -Year 2025 -Month 7 -Day 23 -Hour 21 -Minute 51 -Second 25)
# I78 was text "INATIVO"; it's now a real next-purchase date, so copy the
# date/time number format used elsewhere in the column (matches H78).
$ws.Range("I78").Value = (Get-Date -Year 2026 -Month 1 -Day 23 -Hour 21 -Minute 51 -Second 25)
$ws.Range("I78").NumberFormat = $ws.Range("H78").NumberFormat
$ws.Range("J78").Value = "ATIVO"

# --- Row 116 (id_cliente 28458): daily-purchase counter incremented and
# the client name corrected from "BEMOL S/A" to "RUY MENTA JUNIOR".
$ws.Range("E116").Value = 16928
$ws.Range("H116").Value = (Get-Date -Year 2025 -Month 7 -Day 23 -Hour 15 -Minute 37 -Second 24)
$ws.Range("I116").Value = (Get-Date -Year 2025 -Month 7 -Day 24 -Hour 15 -Minute 37 -Second 24)
$ws.Range("K116").Value = "RUY MENTA JUNIOR"
